$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 827.9524
$ws.Range("J17").Value = 819.35
$ws.Range("L17").Value = 2458.05
$ws.Range("N17").Value = -2794.05
$ws.Range("H62").Value = 4789.3335
$ws.Range("I62").Value = 4628.6924
$ws.Range("J62").Value = 5207
$ws.Range("K62").Value = 4628.6924
$ws.Range("L62").Value = 5207
$ws.Range("M62").Value = -4004.6924
$ws.Range("N62").Value = -6455
$ws.Range("H65").Value = 4789.3335
$ws.Range("I65").Value = 4628.6924
$ws.Range("J65").Value = 5207
$ws.Range("K65").Value = 23143.462
$ws.Range("L65").Value = 26035
$ws.Range("M65").Value = -20023.462
$ws.Range("N65").Value = -32275
$ws.Range("H113").Value = 6292.9644
$ws.Range("I113").Value = 6263.385
$ws.Range("K113").Value = 6263.385
$ws.Range("M113").Value = -3009.385
$ws.Range("H116").Value = 15416.692
$ws.Range("I116").Value = 20326.588
$ws.Range("J116").Value = 6142.4443
$ws.Range("K116").Value = 20326.588
$ws.Range("L116").Value = 6142.4443
$ws.Range("M116").Value = -16884.588
$ws.Range("N116").Value = -13026.4443
$ws.Range("H125").Value = 34412
$ws.Range("I125").Value = 54270.668
$ws.Range("K125").Value = 488436.012
$ws.Range("M125").Value = -485976.012
$ws.Range("H132").Value = 21680.486
$ws.Range("I132").Value = 22712.637
$ws.Range("J132").Value = 4650
$ws.Range("K132").Value = 68137.91099999999
$ws.Range("L132").Value = 13950
$ws.Range("M132").Value = -65607.91099999999
$ws.Range("N132").Value = -19010
$ws.Range("H137").Value = 21577.518
$ws.Range("J137").Value = 30933.1
$ws.Range("L137").Value = 92799.29999999999
$ws.Range("N137").Value = -97899.29999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18260.967
$ws.Range("I32").Value = 20183.81
$ws.Range("K32").Value = 20183.81
$ws.Range("M32").Value = -19896.81
$ws.Range("H45").Value = 3612.6843
$ws.Range("I45").Value = 2059.4167
$ws.Range("K45").Value = 2059.4167
$ws.Range("M45").Value = -1682.4167
$ws.Range("H88").Value = 9903.833000000001
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 9903.833000000001
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 9903.833000000001
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -10715.833
$ws.Range("H91").Value = 9903.833000000001
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 9903.833000000001
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 9903.833000000001
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -12711.833
$ws.Range("H122").Value = 1730.9642
$ws.Range("I122").Value = 1477.25
$ws.Range("K122").Value = 4431.75
$ws.Range("M122").Value = -1981.75
$ws.Range("H127").Value = 101833.336
$ws.Range("J127").Value = 101833.336
$ws.Range("L127").Value = 101833.336
$ws.Range("N127").Value = -111753.336
$ws.Range("H132").Value = 1888
$ws.Range("I132").Value = 1165.1428
$ws.Range("J132").Value = 2900
$ws.Range("K132").Value = 3495.4284
$ws.Range("L132").Value = 8700
$ws.Range("M132").Value = -965.4284000000002
$ws.Range("N132").Value = -13760

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 999.6
$ws.Range("I22").Value = 999.6667
$ws.Range("J22").Value = 999.5
$ws.Range("K22").Value = 999.6667
$ws.Range("L22").Value = 999.5
$ws.Range("M22").Value = -826.6667
$ws.Range("N22").Value = -1345.5
$ws.Range("H94").Value = 7234.625
$ws.Range("I94").Value = 8354.23
$ws.Range("J94").Value = 2383
$ws.Range("K94").Value = 8354.23
$ws.Range("L94").Value = 2383
$ws.Range("M94").Value = -7903.23
$ws.Range("N94").Value = -3285
$ws.Range("H107").Value = 4273.077
$ws.Range("I107").Value = 4248.4
$ws.Range("J107").Value = 4355.3335
$ws.Range("K107").Value = 4248.4
$ws.Range("L107").Value = 4355.3335
$ws.Range("M107").Value = -2328.4
$ws.Range("N107").Value = -8195.333500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4349347
$ws.Range("I31").Value = 6667565.5
$ws.Range("J31").Value = 2687.625
$ws.Range("K31").Value = 6667565.5
$ws.Range("L31").Value = 2687.625
$ws.Range("M31").Value = -6667270.5
$ws.Range("N31").Value = -3277.625
$ws.Range("H34").Value = 4349347
$ws.Range("I34").Value = 6667565.5
$ws.Range("J34").Value = 2687.625
$ws.Range("K34").Value = 6667565.5
$ws.Range("L34").Value = 2687.625
$ws.Range("M34").Value = -6667363.5
$ws.Range("N34").Value = -3091.625
$ws.Range("H86").Value = 46539.117
$ws.Range("I86").Value = 61969.816
$ws.Range("J86").Value = 18249.5
$ws.Range("K86").Value = 61969.816
$ws.Range("L86").Value = 18249.5
$ws.Range("M86").Value = -60846.816
$ws.Range("N86").Value = -20495.5
$ws.Range("H89").Value = 46539.117
$ws.Range("I89").Value = 61969.816
$ws.Range("J89").Value = 18249.5
$ws.Range("K89").Value = 309849.08
$ws.Range("L89").Value = 91247.5
$ws.Range("M89").Value = -304233.08
$ws.Range("N89").Value = -102479.5
$ws.Range("H134").Value = 3294.6191
$ws.Range("I134").Value = 2869.8235
$ws.Range("J134").Value = 5100
$ws.Range("K134").Value = 8609.470499999999
$ws.Range("L134").Value = 15300
$ws.Range("M134").Value = -6074.470499999999
$ws.Range("N134").Value = -20370

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 18733.166
$ws.Range("I29").Value = 800
$ws.Range("J29").Value = 20363.455
$ws.Range("K29").Value = 800
$ws.Range("L29").Value = 20363.455
$ws.Range("M29").Value = -510
$ws.Range("N29").Value = -20943.455
$ws.Range("H80").Value = 11601.111
$ws.Range("I80").Value = 1736.3334
$ws.Range("J80").Value = 16533.5
$ws.Range("K80").Value = 1736.3334
$ws.Range("L80").Value = 16533.5
$ws.Range("M80").Value = -738.3334
$ws.Range("N80").Value = -18529.5
$ws.Range("H83").Value = 11601.111
$ws.Range("I83").Value = 1736.3334
$ws.Range("J83").Value = 16533.5
$ws.Range("K83").Value = 8681.666999999999
$ws.Range("L83").Value = 82667.5
$ws.Range("M83").Value = -3689.666999999999
$ws.Range("N83").Value = -92651.5
$ws.Range("H102").Value = 15218.135
$ws.Range("I102").Value = 20745.076
$ws.Range("K102").Value = 20745.076
$ws.Range("M102").Value = -19123.076
$ws.Range("H113").Value = 3027.1
$ws.Range("I113").Value = 2852.6875
$ws.Range("K113").Value = 2852.6875
$ws.Range("M113").Value = -682.6875
$ws.Range("H132").Value = 1926.2727
$ws.Range("I132").Value = 1625.1428
$ws.Range("J132").Value = 2453.25
$ws.Range("K132").Value = 4875.428400000001
$ws.Range("L132").Value = 7359.75
$ws.Range("M132").Value = -2345.428400000001
$ws.Range("N132").Value = -12419.75
$ws.Range("H134").Value = 62163
$ws.Range("J134").Value = 62163
$ws.Range("L134").Value = 186489
$ws.Range("N134").Value = -191559

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1676.6666
$ws.Range("I22").Value = 1288.8889
$ws.Range("K22").Value = 1288.8889
$ws.Range("M22").Value = -993.8888999999999
$ws.Range("H27").Value = 1676.6666
$ws.Range("I27").Value = 1288.8889
$ws.Range("K27").Value = 1288.8889
$ws.Range("M27").Value = -1181.8889
$ws.Range("H40").Value = 3114.5
$ws.Range("I40").Value = 2988.0715
$ws.Range("K40").Value = 2988.0715
$ws.Range("M40").Value = -2852.0715
$ws.Range("H68").Value = 4047.9375
$ws.Range("I68").Value = 3377
$ws.Range("K68").Value = 3377
$ws.Range("M68").Value = -2628
$ws.Range("H71").Value = 4047.9375
$ws.Range("I71").Value = 3377
$ws.Range("K71").Value = 16885
$ws.Range("M71").Value = -13141
$ws.Range("H76").Value = 14500
$ws.Range("I76").Value = 14500
$ws.Range("K76").Value = 14500
$ws.Range("M76").Value = -14162
$ws.Range("H79").Value = 14500
$ws.Range("I79").Value = 14500
$ws.Range("K79").Value = 14500
$ws.Range("M79").Value = -13330
$ws.Range("H122").Value = 3047.647
$ws.Range("I122").Value = 3013.5557
$ws.Range("J122").Value = 3086
$ws.Range("K122").Value = 9040.667099999999
$ws.Range("L122").Value = 9258
$ws.Range("M122").Value = -6590.667099999999
$ws.Range("N122").Value = -14158
$ws.Range("H132").Value = 2936.375
$ws.Range("I132").Value = 2693.5483
$ws.Range("J132").Value = 3772.7778
$ws.Range("K132").Value = 8080.644899999999
$ws.Range("L132").Value = 11318.3334
$ws.Range("M132").Value = -5550.644899999999
$ws.Range("N132").Value = -16378.3334
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("M133").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 10446959
$ws.Range("I122").Value = 12230088
$ws.Range("K122").Value = 36690264
$ws.Range("M122").Value = -36687814

Write-Output "Applied changes to ALC, ARM, BSM, CRP, GSM, LTW, WVR"